$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update category names (column B) with "Kain " prefix
$ws.Range("B2").Value = "Kain Kaos"
$ws.Range("B3").Value = "Kain Jaket"
$ws.Range("B4").Value = "Kain Jersey"
$ws.Range("B5").Value = "Kain Kemeja"
$ws.Range("B6").Value = "Kain Celana"

# Update slugs (column C) to hyphenated "kain-xxx" form
$ws.Range("C2").Value = "kain-kaos"
$ws.Range("C3").Value = "kain-jaket"
$ws.Range("C4").Value = "kain-jersey"
$ws.Range("C5").Value = "kain-kemeja"
$ws.Range("C6").Value = "kain-celana"

# Remove the last row (previously "Sandal" / "sandal")
$ws.Rows.Item(7).Delete()

# Update selection to reflect the new active cell
$ws.Range("C7").Select()
